# Insert a new slide at the very beginning of the deck using the
# "Content with Caption" layout (ppLayoutContentWithCaption = 35), then
# fill in its placeholders:
#   - Title            -> "Slide with a content and caption"
#   - Content Placeholder (idx 1) -> left empty
#   - Caption Placeholder (body, idx 2, half size) -> "caption text"
#
# All pre-existing slides simply shift down by one position; their
# content is left untouched.

$p = $ppt.ActivePresentation

$s = $p.Slides.Add(1, 35)

$s.Shapes.Item(1).TextFrame.TextRange.Text = "Slide with a content and caption"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "caption text"
